$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 122
$ws.Range("I9").Value = 96
$ws.Range("K9").Value = 96
$ws.Range("M9").Value = 73
$ws.Range("H11").Value = 1264.375
$ws.Range("I11").Value = 1264.375
$ws.Range("K11").Value = 1264.375
$ws.Range("M11").Value = -1124.375
$ws.Range("H33").Value = 98.166664
$ws.Range("I33").Value = 100.72727
$ws.Range("J33").Value = 70
$ws.Range("K33").Value = 100.72727
$ws.Range("L33").Value = 70
$ws.Range("M33").Value = 128.27273
$ws.Range("N33").Value = -528
$ws.Range("H38").Value = 8
$ws.Range("I38").Value = 8
$ws.Range("K38").Value = 24
$ws.Range("M38").Value = 348
$ws.Range("H42").Value = 59.57143
$ws.Range("I42").Value = 66.75
$ws.Range("K42").Value = 200.25
$ws.Range("M42").Value = 29.75
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = $null
$ws.Range("H138").Value = 8124.75
$ws.Range("J138").Value = 9750
$ws.Range("L138").Value = 29250
$ws.Range("N138").Value = -39530
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1680
$ws.Range("I74").Value = 1680
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1680
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -806
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 1680
$ws.Range("I77").Value = 1680
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8400
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4032
$ws.Range("N77").Value = $null
$ws.Range("H92").Value = 90199.8
$ws.Range("I92").Value = 63000
$ws.Range("J92").Value = 96999.75
$ws.Range("K92").Value = 63000
$ws.Range("L92").Value = 96999.75
$ws.Range("M92").Value = -60504
$ws.Range("N92").Value = -101991.75
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null
$ws.Range("H97").Value = 1755
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null
$ws.Range("H102").Value = 1253
$ws.Range("I102").Value = 1253
$ws.Range("K102").Value = 1253
$ws.Range("M102").Value = 369
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 1059
$ws.Range("I41").Value = 1059
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1059
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -631
$ws.Range("N41").Value = $null
$ws.Range("H68").Value = 49800
$ws.Range("J68").Value = 49800
$ws.Range("L68").Value = 49800
$ws.Range("N68").Value = -51298
$ws.Range("H69").Value = 43000.6
$ws.Range("H71").Value = 49800
$ws.Range("J71").Value = 49800
$ws.Range("L71").Value = 149400
$ws.Range("N71").Value = -156888
$ws.Range("H72").Value = 43000.6
$ws.Range("H74").Value = 26599.4
$ws.Range("J74").Value = 29999.25
$ws.Range("L74").Value = 29999.25
$ws.Range("N74").Value = -31747.25
$ws.Range("H77").Value = 26599.4
$ws.Range("J77").Value = 29999.25
$ws.Range("L77").Value = 89997.75
$ws.Range("N77").Value = -98733.75
$ws.Range("H107").Value = 594.3
$ws.Range("I107").Value = 630.875
$ws.Range("J107").Value = 448
$ws.Range("K107").Value = 630.875
$ws.Range("L107").Value = 448
$ws.Range("M107").Value = 1289.125
$ws.Range("N107").Value = -4288
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 224.08696
$ws.Range("I2").Value = 241.38095
$ws.Range("K2").Value = 1448.2857
$ws.Range("M2").Value = -1335.2857
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("H109").Value = 5995
$ws.Range("I109").Value = 5995
$ws.Range("K109").Value = 17985
$ws.Range("M109").Value = -16945
$ws.Range("H114").Value = 2313.1667
$ws.Range("J114").Value = 2862.75
$ws.Range("L114").Value = 8588.25
$ws.Range("N114").Value = -15096.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 6500
$ws.Range("I33").Value = 6500
$ws.Range("K33").Value = 6500
$ws.Range("M33").Value = -6248
$ws.Range("H36").Value = 10000
$ws.Range("J36").Value = 10000
$ws.Range("L36").Value = 10000
$ws.Range("N36").Value = -10970
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = $null
$ws.Range("H122").Value = 4097.875
$ws.Range("I122").Value = 3572.2856
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 10716.8568
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -8266.856800000001
$ws.Range("N122").Value = -28231
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 255
$ws.Range("J31").Value = 495
$ws.Range("L31").Value = 495
$ws.Range("N31").Value = -991
$ws.Range("H50").Value = 15000
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").Value = $null
$ws.Range("H63").Value = 25000
$ws.Range("J63").Value = 25000
$ws.Range("L63").Value = 25000
$ws.Range("N63").Value = -26498
$ws.Range("H66").Value = 25000
$ws.Range("J66").Value = 25000
$ws.Range("L66").Value = 75000
$ws.Range("N66").Value = -82488
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null
$ws.Range("H100").Value = 5118
$ws.Range("I100").Value = 5118
$ws.Range("K100").Value = 5118
$ws.Range("M100").Value = -4577
$ws.Range("H132").Value = 9213.429
$ws.Range("I132").Value = 4899
$ws.Range("J132").Value = 19999.5
$ws.Range("K132").Value = 14697
$ws.Range("L132").Value = 59998.5
$ws.Range("M132").Value = -12167
$ws.Range("N132").Value = -65058.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null
$ws.Range("H51").Value = 3000
$ws.Range("I51").Value = 3000
$ws.Range("K51").Value = 3000
$ws.Range("M51").Value = -2490
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = $null
$ws.Range("H64").Value = 40000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null
$ws.Range("H67").Value = 40000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null
$ws.Range("H70").Value = 25000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = $null
$ws.Range("H73").Value = 25000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = $null
$ws.Range("H75").Value = 25000
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = $null
$ws.Range("H78").Value = 25000
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = $null
$ws.Range("H94").Value = 35331
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 35331
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 35331
$ws.Range("M94").Value = $null
$ws.Range("N94").Value = -37133
$ws.Range("H103").Value = 44000
$ws.Range("J103").Value = 44000
$ws.Range("L103").Value = 44000
$ws.Range("N103").Value = -46344
